# #5: property boat&car done
# Update the "汽車" (car) sheet (sheet3 / 3rd worksheet tab) so that it
# follows the same schema as the other property sheets: a proper header
# row with field names, and the data row extended with the common
# trailing metadata columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index). A new field
# "capacity" (replacing the unlabeled engine-displacement column) is
# introduced as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Header row (row 1): turn it into real field names ---------------
$ws.Range("B1").Value2 = "name"
$ws.Range("C1").Value2 = "capacity"
$ws.Range("D1").Value2 = "owner"
$ws.Range("E1").Value2 = "register_date"
$ws.Range("F1").Value2 = "register_reason"
$ws.Range("G1").Value2 = "acquire_value"
$ws.Range("H1").Value2 = "property_category"
$ws.Range("I1").Value2 = "category"
$ws.Range("J1").Value2 = "date"
$ws.Range("K1").Value2 = "legislator_name"
$ws.Range("L1").Value2 = "legislator_id"
$ws.Range("M1").Value2 = "source_file"
$ws.Range("N1").Value2 = "index"

# --- Data row (row 2): keep existing values, append new columns ------
# (A2, C2, D2, F2 already hold the correct values: 40 / 3498 / 曾月桂 / 買賣)
$ws.Range("B2").Value2 = "BENZE350"
$ws.Range("E2").Value2 = "101年01月20曰"
$ws.Range("G2").Value2 = "1350000(2006年3月出廠）"
$ws.Range("H2").Value2 = "land"
$ws.Range("I2").Value2 = "normal"
$ws.Range("J2").Value2 = "2012-03-14"
$ws.Range("K2").Value2 = "林德福"
$ws.Range("L2").Value2 = 908
$ws.Range("M2").Value2 = "tmp82d01"
$ws.Range("N2").Value2 = 40
